# Update NATMI LR-pair output (Efna1-Epha4) with newly recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 23.630375
$ws.Range("H2").Value = 70.89112499999999
$ws.Range("I2").Value = 0.9002398112414131
$ws.Range("J2").Value = 0.9002398112414129
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.081040666666667
$ws.Range("N2").Value = 24.243122
$ws.Range("O2").Value = 0.4661250698616886
$ws.Range("P2").Value = 0.4661250698616886
$ws.Range("Q2").Value = 190.9580213435833
$ws.Range("R2").Value = 1718.62219209225
$ws.Range("S2").Value = 0.419624344907177
$ws.Range("T2").Value = 0.419624344907177

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 23.630375
$ws.Range("H3").Value = 70.89112499999999
$ws.Range("I3").Value = 0.9002398112414131
$ws.Range("J3").Value = 0.9002398112414129
$ws.Range("O3").Value = 0.4037865631294714
$ws.Range("P3").Value = 0.4037865631294715
$ws.Range("Q3").Value = 165.4197298660833
$ws.Range("R3").Value = 1488.77756879475
$ws.Range("S3").Value = 0.3635047393734943
$ws.Range("T3").Value = 0.3635047393734943

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 23.630375
$ws.Range("H4").Value = 70.89112499999999
$ws.Range("I4").Value = 0.9002398112414131
$ws.Range("J4").Value = 0.9002398112414129
$ws.Range("O4").Value = 0.1300883670088399
$ws.Range("P4").Value = 0.1300883670088399
$ws.Range("Q4").Value = 53.29345870883332
$ws.Range("R4").Value = 479.6411283794999
$ws.Range("S4").Value = 0.1171107269607417
$ws.Range("T4").Value = 0.1171107269607417

# Row 5: FAPs -> ECs
$ws.Range("I5").Value = 0.06214870537054815
$ws.Range("J5").Value = 0.06214870537054815
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.081040666666667
$ws.Range("N5").Value = 24.243122
$ws.Range("O5").Value = 0.4661250698616886
$ws.Range("P5").Value = 0.4661250698616886
$ws.Range("Q5").Value = 13.18292488116
$ws.Range("R5").Value = 118.64632393044
$ws.Range("S5").Value = 0.02896906963266026
$ws.Range("T5").Value = 0.02896906963266026

# Row 6: FAPs -> FAPs
$ws.Range("I6").Value = 0.06214870537054815
$ws.Range("J6").Value = 0.06214870537054815
$ws.Range("O6").Value = 0.4037865631294714
$ws.Range("P6").Value = 0.4037865631294715
$ws.Range("S6").Value = 0.02509481214451976
$ws.Range("T6").Value = 0.02509481214451976

# Row 7: FAPs -> MuSCs
$ws.Range("I7").Value = 0.06214870537054815
$ws.Range("J7").Value = 0.06214870537054815
$ws.Range("O7").Value = 0.1300883670088399
$ws.Range("P7").Value = 0.1300883670088399
$ws.Range("S7").Value = 0.008084823593368127
$ws.Range("T7").Value = 0.008084823593368127

# Row 8: MuSCs -> ECs
$ws.Range("I8").Value = 0.03761148338803896
$ws.Range("J8").Value = 0.03761148338803896
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.081040666666667
$ws.Range("N8").Value = 24.243122
$ws.Range("O8").Value = 0.4661250698616886
$ws.Range("P8").Value = 0.4661250698616886
$ws.Range("Q8").Value = 7.978112451695333
$ws.Range("R8").Value = 71.803012065258
$ws.Range("S8").Value = 0.0175316553218514
$ws.Range("T8").Value = 0.0175316553218514

# Row 9: MuSCs -> FAPs
$ws.Range("I9").Value = 0.03761148338803896
$ws.Range("J9").Value = 0.03761148338803896
$ws.Range("O9").Value = 0.4037865631294714
$ws.Range("P9").Value = 0.4037865631294715
$ws.Range("S9").Value = 0.01518701161145746
$ws.Range("T9").Value = 0.01518701161145746

# Row 10: MuSCs -> MuSCs
$ws.Range("I10").Value = 0.03761148338803896
$ws.Range("J10").Value = 0.03761148338803896
$ws.Range("O10").Value = 0.1300883670088399
$ws.Range("P10").Value = 0.1300883670088399
$ws.Range("S10").Value = 0.004892816454730099
$ws.Range("T10").Value = 0.004892816454730099
